$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B200").HorizontalAlignment = -4108
$ws.Range("B200").VerticalAlignment = -4108
